{"js": "// Update each \"three-digit \u00f7 one-digit\" answer cell to the regenerated value.\n// `replacements` is the literal old-text -> new-text mapping from the target diff;\n// every cell's text is unique in the document, so a plain text search locates it.\nconst replacements = [\n  [\"168\u00f77=24, 0\", \"127\u00f72=63, 1\"],\n  [\"195\u00f78=24, 3\", \"576\u00f72=288, 0\"],\n  [\"140\u00f73=46, 2\", \"501\u00f78=62, 5\"],\n  [\"658\u00f76=109, 4\", \"764\u00f76=127, 2\"],\n  [\"849\u00f75=169, 4\", \"326\u00f75=65, 1\"],\n  [\"803\u00f79=89, 2\", \"638\u00f74=159, 2\"],\n  [\"616\u00f72=308, 0\", \"126\u00f73=42, 0\"],\n  [\"961\u00f75=192, 1\", \"892\u00f78=111, 4\"],\n  [\"554\u00f79=61, 5\", \"107\u00f78=13, 3\"],\n  [\"330\u00f75=66, 0\", \"759\u00f75=151, 4\"],\n  [\"946\u00f77=135, 1\", \"245\u00f77=35, 0\"],\n  [\"592\u00f76=98, 4\", \"127\u00f75=25, 2\"],\n  [\"156\u00f77=22, 2\", \"547\u00f74=136, 3\"],\n  [\"338\u00f76=56, 2\", \"819\u00f72=409, 1\"],\n  [\"894\u00f79=99, 3\", \"499\u00f73=166, 1\"],\n  [\"142\u00f78=17, 6\", \"702\u00f74=175, 2\"],\n  [\"395\u00f79=43, 8\", \"807\u00f75=161, 2\"],\n  [\"552\u00f79=61, 3\", \"509\u00f78=63, 5\"],\n  [\"832\u00f75=166, 2\", \"955\u00f73=318, 1\"],\n  [\"127\u00f76=21, 1\", \"815\u00f73=271, 2\"],\n  [\"644\u00f74=161, 0\", \"197\u00f79=21, 8\"],\n  [\"788\u00f75=157, 3\", \"307\u00f76=51, 1\"],\n  [\"887\u00f72=443, 1\", \"572\u00f76=95, 2\"],\n  [\"648\u00f73=216, 0\", \"646\u00f73=215, 1\"],\n  [\"675\u00f76=112, 3\", \"519\u00f74=129, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-division answer cell with its new value.\n# The pairs below come 1:1 from the authoritative diff (old text -> new text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"168\u00f77=24, 0\", \"127\u00f72=63, 1\"),\n    @(\"195\u00f78=24, 3\", \"576\u00f72=288, 0\"),\n    @(\"140\u00f73=46, 2\", \"501\u00f78=62, 5\"),\n    @(\"658\u00f76=109, 4\", \"764\u00f76=127, 2\"),\n    @(\"849\u00f75=169, 4\", \"326\u00f75=65, 1\"),\n    @(\"803\u00f79=89, 2\", \"638\u00f74=159, 2\"),\n    @(\"616\u00f72=308, 0\", \"126\u00f73=42, 0\"),\n    @(\"961\u00f75=192, 1\", \"892\u00f78=111, 4\"),\n    @(\"554\u00f79=61, 5\", \"107\u00f78=13, 3\"),\n    @(\"330\u00f75=66, 0\", \"759\u00f75=151, 4\"),\n    @(\"946\u00f77=135, 1\", \"245\u00f77=35, 0\"),\n    @(\"592\u00f76=98, 4\", \"127\u00f75=25, 2\"),\n    @(\"156\u00f77=22, 2\", \"547\u00f74=136, 3\"),\n    @(\"338\u00f76=56, 2\", \"819\u00f72=409, 1\"),\n    @(\"894\u00f79=99, 3\", \"499\u00f73=166, 1\"),\n    @(\"142\u00f78=17, 6\", \"702\u00f74=175, 2\"),\n    @(\"395\u00f79=43, 8\", \"807\u00f75=161, 2\"),\n    @(\"552\u00f79=61, 3\", \"509\u00f78=63, 5\"),\n    @(\"832\u00f75=166, 2\", \"955\u00f73=318, 1\"),\n    @(\"127\u00f76=21, 1\", \"815\u00f73=271, 2\"),\n    @(\"644\u00f74=161, 0\", \"197\u00f79=21, 8\"),\n    @(\"788\u00f75=157, 3\", \"307\u00f76=51, 1\"),\n    @(\"887\u00f72=443, 1\", \"572\u00f76=95, 2\"),\n    @(\"648\u00f73=216, 0\", \"646\u00f73=215, 1\"),\n    @(\"675\u00f76=112, 3\", \"519\u00f74=129, 3\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 -> wdFindContinue (search whole story), Replace:=2 -> wdReplaceAll.\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Host \"WARNING: no match found for: $oldText\"\n    }\n}\n\n"}
